$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=0.366566; H=1.099698; I=0.5689653834353526; J=0.5689653834353527; K=3; M=4.166450999999999; N=12.499353; O=0.7126954333415383; P=0.7126954333415383; Q=1.527279277266; R=13.745513495394; S=0.4054990305037932; T=0.4054990305037932 }
    3 = @{ E=3; G=0.366566; H=1.099698; I=0.5689653834353526; J=0.5689653834353527; K=3; M=0.6655859999999999; N=1.996758; O=0.1138523176430159; P=0.1138523176430159; Q=0.243981197676; R=2.195830779084; S=0.06477802756276207; T=0.06477802756276209 }
    4 = @{ E=3; G=0.366566; H=1.099698; I=0.5689653834353526; J=0.5689653834353527; K=3; M=1.01401; N=3.04203; O=0.1734522490154458; P=0.1734522490154458; Q=0.3717015896600001; R=3.34531430694; S=0.09868832536879742; T=0.09868832536879742 }
    5 = @{ E=3; G=0.2777016666666667; H=0.833105; I=0.4310346165646473; J=0.4310346165646473; K=3; M=4.166450999999999; N=12.499353; O=0.7126954333415383; P=0.7126954333415383; Q=1.157030386785; R=10.413273481065; S=0.3071964028377451; T=0.3071964028377451 }
    6 = @{ E=3; G=0.2777016666666667; H=0.833105; I=0.4310346165646473; J=0.4310346165646473; K=3; M=0.6655859999999999; N=1.996758; O=0.1138523176430159; P=0.1138523176430159; Q=0.18483434151; R=1.66350907359; S=0.04907429008025377; T=0.04907429008025377 }
    7 = @{ E=3; G=0.2777016666666667; H=0.833105; I=0.4310346165646473; J=0.4310346165646473; K=3; M=1.01401; N=3.04203; O=0.1734522490154458; P=0.1734522490154458; Q=0.2815922670166667; R=2.53433040315; S=0.07476392364664843; T=0.07476392364664841 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
